$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (technologische Museum Sr. Maj. des Kaisers) gains full data,
#     copied from row 3's coordinate/source columns ---
$ws.Range("C4").Value2 = "Jn Gebäude des voltechi=schen Jnstituts 1. Hof, Direk=tionstiege, 1. Stock="
$ws.Range("E4").Value2 = 16.369911999999999
$ws.Range("F4").Value2 = 48.198887999999997
$ws.Range("G4").Value2 = "Wikipedia, GeoHack"
$ws.Range("H4").Value2 = "https://geohack.toolforge.org/geohack.php?pagename=K.k._Polytechnisches_Institut&language=de&params=48.198888_N_16.369912_E_region:AT-9_type:landmark"
$ws.Range("E4").Style = $ws.Range("E3").Style

# --- New addresses for rows 5, 6, 7 ---
$ws.Range("C5").Value2 = "04., Taubstummengasse 13-17"
$ws.Range("C6").Value2 = "04., Favoritenstraße 15"
$ws.Range("C7").Value2 = "04., Favoritenstraße 16"

# --- New source URLs for rows 6, 7, 5 ---
$ws.Range("H6").Value2 = "https://www.geschichtewiki.wien.gv.at/Theresianische_Akademie_(Geb%C3%A4ude)"
$ws.Range("H7").Value2 = "https://www.geschichtewiki.wien.gv.at/Theresianische_Akademie_(Geb%C3%A4ude)"
$ws.Range("H5").Value2 = "https://www.geschichtewiki.wien.gv.at/Taubstummeninstitut"

# --- Rename header columns Longitude_x / Latitude_x ---
$ws.Range("E1").Value2 = "Longitude"
$ws.Range("F1").Value2 = "Latitude"

# --- New "Quelle" for rows 8 and 9 ---
$ws.Range("G8").Value2 = "QGIS - Estimation"
$ws.Range("G9").Value2 = "QGIS - Estimation"

# --- New "Sicherheit" values for rows 5, 6, 7 ---
$ws.Range("D5").Value2 = "mittel"
$ws.Range("D6").Value2 = "mittel"
$ws.Range("D7").Value2 = "mittel"

# --- Remaining (no new shared strings) updates ---
$ws.Range("D4").Value2 = "hoch"
$ws.Range("D8").Value2 = "hoch"
$ws.Range("D9").Value2 = "hoch"

$ws.Range("E5").Value2 = 16.3709411171927
$ws.Range("F5").Value2 = 48.195140442994003
$ws.Range("G5").Value2 = "Wien Geschichte Wiki, Google Maps"

$ws.Range("E6").Value2 = 16.3714736393609
$ws.Range("F6").Value2 = 48.193537549633099
$ws.Range("G6").Value2 = "Wien Geschichte Wiki, Google Maps"

$ws.Range("E7").Value2 = 16.3714736393609
$ws.Range("F7").Value2 = 48.193537549633099
$ws.Range("G7").Value2 = "Wien Geschichte Wiki, Google Maps"

$ws.Range("E8").Value2 = 16.379925897949999
$ws.Range("F8").Value2 = 48.185791939397014

$ws.Range("E9").Value2 = 16.374092267197891
$ws.Range("F9").Value2 = 48.182660301764052

$ws.Range("C9").ClearContents()
$ws.Range("H9").ClearContents()

# --- Stray formatted (empty) cell that appears below the table, a leftover
#     from the editing session, bold-styled like the header font ---
$ws.Range("F14").Font.Bold = $true

# --- Final cursor/selection position left by the editing session ---
[void]$ws.Range("E13").Select()
